$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cycle [6, 8, 7]
$ws.Range("A6:B6").Copy($ws.Range("A1000"))
$ws.Range("D6:N6").Copy($ws.Range("D1000"))
$ws.Range("Q6:R6").Copy($ws.Range("Q1000"))
$ws.Range("AF1000").ClearContents()

$ws.Range("A8:B8").Copy($ws.Range("A6"))
$ws.Range("D8:N8").Copy($ws.Range("D6"))
$ws.Range("Q8:R8").Copy($ws.Range("Q6"))
$ws.Range("AF8").Copy($ws.Range("AF6"))

$ws.Range("A7:B7").Copy($ws.Range("A8"))
$ws.Range("D7:N7").Copy($ws.Range("D8"))
$ws.Range("Q7:R7").Copy($ws.Range("Q8"))
$ws.Range("AF8").ClearContents()

$ws.Range("A1000:B1000").Copy($ws.Range("A7"))
$ws.Range("D1000:N1000").Copy($ws.Range("D7"))
$ws.Range("Q1000:R1000").Copy($ws.Range("Q7"))
$ws.Range("AF7").ClearContents()

$ws.Range("A1000:AF1000").ClearContents()

# Cycle [20, 22]
$ws.Range("A20:B20").Copy($ws.Range("A1001"))
$ws.Range("D20:N20").Copy($ws.Range("D1001"))
$ws.Range("Q20:R20").Copy($ws.Range("Q1001"))
$ws.Range("AF1001").ClearContents()

$ws.Range("A22:B22").Copy($ws.Range("A20"))
$ws.Range("D22:N22").Copy($ws.Range("D20"))
$ws.Range("Q22:R22").Copy($ws.Range("Q20"))
$ws.Range("AF22").Copy($ws.Range("AF20"))

$ws.Range("A1001:B1001").Copy($ws.Range("A22"))
$ws.Range("D1001:N1001").Copy($ws.Range("D22"))
$ws.Range("Q1001:R1001").Copy($ws.Range("Q22"))
$ws.Range("AF22").ClearContents()

$ws.Range("A1001:AF1001").ClearContents()

# Cycle [29, 30]
$ws.Range("A29:B29").Copy($ws.Range("A1002"))
$ws.Range("D29:N29").Copy($ws.Range("D1002"))
$ws.Range("Q29:R29").Copy($ws.Range("Q1002"))
$ws.Range("AF1002").ClearContents()

$ws.Range("A30:B30").Copy($ws.Range("A29"))
$ws.Range("D30:N30").Copy($ws.Range("D29"))
$ws.Range("Q30:R30").Copy($ws.Range("Q29"))
$ws.Range("AF30").Copy($ws.Range("AF29"))

$ws.Range("A1002:B1002").Copy($ws.Range("A30"))
$ws.Range("D1002:N1002").Copy($ws.Range("D30"))
$ws.Range("Q1002:R1002").Copy($ws.Range("Q30"))
$ws.Range("AF30").ClearContents()

$ws.Range("A1002:AF1002").ClearContents()

# Cycle [36, 39, 38, 37]
$ws.Range("A36:B36").Copy($ws.Range("A1003"))
$ws.Range("D36:N36").Copy($ws.Range("D1003"))
$ws.Range("Q36:R36").Copy($ws.Range("Q1003"))
$ws.Range("AF1003").ClearContents()

$ws.Range("A39:B39").Copy($ws.Range("A36"))
$ws.Range("D39:N39").Copy($ws.Range("D36"))
$ws.Range("Q39:R39").Copy($ws.Range("Q36"))
$ws.Range("AF39").Copy($ws.Range("AF36"))

$ws.Range("A38:B38").Copy($ws.Range("A39"))
$ws.Range("D38:N38").Copy($ws.Range("D39"))
$ws.Range("Q38:R38").Copy($ws.Range("Q39"))
$ws.Range("AF39").ClearContents()

$ws.Range("A37:B37").Copy($ws.Range("A38"))
$ws.Range("D37:N37").Copy($ws.Range("D38"))
$ws.Range("Q37:R37").Copy($ws.Range("Q38"))
$ws.Range("AF38").ClearContents()

$ws.Range("A1003:B1003").Copy($ws.Range("A37"))
$ws.Range("D1003:N1003").Copy($ws.Range("D37"))
$ws.Range("Q1003:R1003").Copy($ws.Range("Q37"))
$ws.Range("AF37").ClearContents()

$ws.Range("A1003:AF1003").ClearContents()

# Cycle [41, 43, 42]
$ws.Range("A41:B41").Copy($ws.Range("A1004"))
$ws.Range("D41:N41").Copy($ws.Range("D1004"))
$ws.Range("Q41:R41").Copy($ws.Range("Q1004"))
$ws.Range("AF1004").ClearContents()

$ws.Range("A43:B43").Copy($ws.Range("A41"))
$ws.Range("D43:N43").Copy($ws.Range("D41"))
$ws.Range("Q43:R43").Copy($ws.Range("Q41"))
$ws.Range("AF41").ClearContents()

$ws.Range("A42:B42").Copy($ws.Range("A43"))
$ws.Range("D42:N42").Copy($ws.Range("D43"))
$ws.Range("Q42:R42").Copy($ws.Range("Q43"))
$ws.Range("AF43").ClearContents()

$ws.Range("A1004:B1004").Copy($ws.Range("A42"))
$ws.Range("D1004:N1004").Copy($ws.Range("D42"))
$ws.Range("Q1004:R1004").Copy($ws.Range("Q42"))
$ws.Range("AF42").ClearContents()

$ws.Range("A1004:AF1004").ClearContents()
